$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$theme = $sm.Theme
Write-Host "theme: $theme"
Write-Host "ThemeVariants.Count: $($theme.ThemeVariants.Count)"

$nm = $p.NotesMaster
$theme2 = $nm.Theme
Write-Host "theme2(notesmaster): $theme2"
if ($theme2 -eq $null) { Write-Host "notesmaster theme null" }
